# Feria Lagunitas de Puerto Montt - Haba: add a new weekly price record.
# A new row is inserted at row 30 (pushing the existing rows 30-74 down to
# 31-75) and populated with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 30; this shifts rows 30:74 down to 31:75
# and carries the row-above formatting (incl. the date style on column D).
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new data point.
$ws.Cells.Item(30, 1).Value  = 4
$ws.Cells.Item(30, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(30, 3).Value  = "Los Lagos"
$ws.Cells.Item(30, 4).Value  = 44571
$ws.Cells.Item(30, 5).Value  = 10
$ws.Cells.Item(30, 6).Value  = 100112026
$ws.Cells.Item(30, 7).Value  = "Haba"
$ws.Cells.Item(30, 8).Value  = "Sin especificar"
$ws.Cells.Item(30, 9).Value  = "Primera"
$ws.Cells.Item(30, 10).Value = 80
$ws.Cells.Item(30, 11).Value = 23000
$ws.Cells.Item(30, 12).Value = 23000
$ws.Cells.Item(30, 13).Value = 23000
$ws.Cells.Item(30, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(30, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(30, 16).Value = 920
$ws.Cells.Item(30, 17).Value = 25
$ws.Cells.Item(30, 18).Value = "Hortaliza"
